$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" - updated aggregate metrics
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(3, 2).Value = 1199.66   # Current Capital
$wsSummary.Cells.Item(4, 2).Value = -0.35     # Total P&L $
$wsSummary.Cells.Item(5, 2).Value = -0.08     # Total P&L %
$wsSummary.Cells.Item(6, 2).Value = 91        # Total Trades
$wsSummary.Cells.Item(8, 2).Value = 48        # Losing Trades
$wsSummary.Cells.Item(9, 2).Value = 34.07     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status" - MarketMaking row updated
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Cells.Item(4, 3).Value = 99.66   # Capital
$wsStatus.Cells.Item(4, 4).Value = 91      # Trades
$wsStatus.Cells.Item(4, 5).Value = -0.35   # P&L $
$wsStatus.Cells.Item(4, 6).Value = -0.34   # P&L %
$wsStatus.Cells.Item(4, 7).Value = 34.07   # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade #91 as row 92 to a trade-log sheet while
# keeping the date/time columns as literal text (matching the existing
# inline-string cells instead of being auto-converted to date serials).
# ---------------------------------------------------------------------------
function Add-Trade91Row {
    param($ws)

    $row = 92

    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $timeCell = $ws.Cells.Item($row, 3)
    $timeCell.NumberFormat = "@"
    $timeCell.Value = "15:54:13"
    $timeCell.Style = "Normal"

    $ws.Cells.Item($row, 1).Value = 91
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.79
    $ws.Cells.Item($row, 7).Value = 0.73
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -7.5949
    $ws.Cells.Item($row, 10).Value = -0.06
    $ws.Cells.Item($row, 11).Value = 99.66
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

# ---------------------------------------------------------------------------
# Sheet "All Trades" - append trade #91
# ---------------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade91Row $wsAllTrades

# ---------------------------------------------------------------------------
# Sheet "MarketMaking" - append trade #91 (mirrors "All Trades")
# ---------------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade91Row $wsMarketMaking

Write-Host "Applied trade #91 update across Summary, Strategy Status, All Trades, MarketMaking sheets."
